$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5243.5
$ws.Range("J32").Value = 4665.3076
$ws.Range("L32").Value = 4665.3076
$ws.Range("N32").Value = -5317.3076
$ws.Range("H40").Value = 6866.125
$ws.Range("J40").Value = 8166.6665
$ws.Range("L40").Value = 8166.6665
$ws.Range("N40").Value = -8516.666499999999
$ws.Range("H43").Value = 5473.875
$ws.Range("J43").Value = 6722.75
$ws.Range("L43").Value = 6722.75
$ws.Range("N43").Value = -6860.75
$ws.Range("H51").Value = 9210.210999999999
$ws.Range("J51").Value = 9533
$ws.Range("L51").Value = 9533
$ws.Range("N51").Value = -10501
$ws.Range("H55").Value = 385.16666
$ws.Range("I55").Value = 69
$ws.Range("J55").Value = 543.25
$ws.Range("K55").Value = 69
$ws.Range("L55").Value = 543.25
$ws.Range("M55").Value = 145
$ws.Range("N55").Value = -971.25
$ws.Range("H80").Value = 466.16666
$ws.Range("I80").Value = 492.14285
$ws.Range("J80").Value = 429.8
$ws.Range("K80").Value = 1476.42855
$ws.Range("L80").Value = 1289.4
$ws.Range("M80").Value = -478.4285500000001
$ws.Range("N80").Value = -3285.4
$ws.Range("H83").Value = 466.16666
$ws.Range("I83").Value = 492.14285
$ws.Range("J83").Value = 429.8
$ws.Range("K83").Value = 4429.28565
$ws.Range("L83").Value = 3868.2
$ws.Range("M83").Value = 562.7143500000002
$ws.Range("N83").Value = -13852.2
$ws.Range("H92").Value = 187.08333
$ws.Range("I92").Value = 187.08333
$ws.Range("K92").Value = 187.08333
$ws.Range("M92").Value = 1060.91667
$ws.Range("H98").Value = 2110
$ws.Range("I98").Value = 2110
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2110
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -612
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 2110
$ws.Range("I122").Value = 2110
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6330
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3880
$ws.Range("N122").Value = -3880
$ws.Range("H129").Value = 1705.5
$ws.Range("I129").Value = 1447.25
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 4341.75
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 658.25
$ws.Range("N129").Value = -16666
$ws.Range("H138").Value = 41669024
$ws.Range("I138").Value = 3054.7
$ws.Range("J138").Value = 71430424
$ws.Range("K138").Value = 9164.099999999999
$ws.Range("L138").Value = 214291272
$ws.Range("M138").Value = -4024.099999999999
$ws.Range("N138").Value = -214301552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7008.273
$ws.Range("J46").Value = 7243.4443
$ws.Range("L46").Value = 7243.4443
$ws.Range("N46").Value = -7881.4443
$ws.Range("H61").Value = 12201107
$ws.Range("I61").Value = 14710840
$ws.Range("J61").Value = 10976.143
$ws.Range("K61").Value = 14710840
$ws.Range("L61").Value = 10976.143
$ws.Range("M61").Value = -14710628
$ws.Range("N61").Value = -11400.143
$ws.Range("H136").Value = 12201107
$ws.Range("I136").Value = 14710840
$ws.Range("J136").Value = 10976.143
$ws.Range("K136").Value = 44132520
$ws.Range("L136").Value = 32928.429
$ws.Range("M136").Value = -44129970
$ws.Range("N136").Value = -38028.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4108.1514
$ws.Range("I134").Value = 4283.357
$ws.Range("K134").Value = 12850.071
$ws.Range("M134").Value = -10315.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J22").Value = 762.25
$ws.Range("L22").Value = 762.25
$ws.Range("N22").Value = -1462.25
$ws.Range("H68").Value = 41281.75
$ws.Range("J68").Value = 41281.75
$ws.Range("L68").Value = 41281.75
$ws.Range("N68").Value = -42779.75
$ws.Range("H71").Value = 41281.75
$ws.Range("J71").Value = 41281.75
$ws.Range("L71").Value = 123845.25
$ws.Range("N71").Value = -131333.25
$ws.Range("H74").Value = 36523.145
$ws.Range("J74").Value = 36523.145
$ws.Range("L74").Value = 36523.145
$ws.Range("N74").Value = -38271.145
$ws.Range("H77").Value = 36523.145
$ws.Range("J77").Value = 36523.145
$ws.Range("L77").Value = 109569.435
$ws.Range("N77").Value = -118305.435
$ws.Range("H102").Value = 195000
$ws.Range("J102").Value = 360000
$ws.Range("L102").Value = 360000
$ws.Range("N102").Value = -364868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1333
$ws.Range("I22").Value = 1333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3999
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3830
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1333
$ws.Range("I27").Value = 1333
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3999
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3897
$ws.Range("N27").ClearContents()
$ws.Range("H128").Value = 516996.5
$ws.Range("I128").Value = 516996.5
$ws.Range("K128").Value = 1550989.5
$ws.Range("M128").Value = -1546009.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 34959
$ws.Range("J59").Value = 34959
$ws.Range("L59").Value = 34959
$ws.Range("N59").Value = -36125
$ws.Range("H113").Value = 3069
$ws.Range("I113").Value = 3096.6
$ws.Range("K113").Value = 3096.6
$ws.Range("M113").Value = -926.5999999999999
$ws.Range("H126").Value = 3956.0715
$ws.Range("I126").Value = 3098.3333
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 9294.999899999999
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -6824.999899999999
$ws.Range("N126").Value = -21440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 610
$ws.Range("I55").Value = 792.375
$ws.Range("J55").Value = 427.625
$ws.Range("K55").Value = 792.375
$ws.Range("L55").Value = 427.625
$ws.Range("M55").Value = -619.375
$ws.Range("N55").Value = -773.625
$ws.Range("H99").Value = 46333
$ws.Range("I99").Value = 46333
$ws.Range("K99").Value = 46333
$ws.Range("M99").Value = -43338
$ws.Range("H139").Value = 62550
$ws.Range("J139").Value = 85100
$ws.Range("L139").Value = 85100
$ws.Range("N139").Value = -95380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11931.143
$ws.Range("I41").Value = 6749.5
$ws.Range("J41").Value = 12794.75
$ws.Range("K41").Value = 6749.5
$ws.Range("L41").Value = 12794.75
$ws.Range("M41").Value = -6359.5
$ws.Range("N41").Value = -13574.75
$ws.Range("H54").Value = 28800
$ws.Range("J54").Value = 28800
$ws.Range("L54").Value = 28800
$ws.Range("N54").Value = -29840
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
